# Update loading_percent values on Sheet1 for the "380 kV" case run.
# Only columns B,D,E,F,G,H,I,J,L,M,N change for rows 2-25 (columns C,K,O stay 0,
# and column A holds the unchanged row index). Values below are the new results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16.78387890987024
$ws.Range("D2").Value = 8.351231686943816
$ws.Range("E2").Value = 13.81062963640493
$ws.Range("F2").Value = 37.64517756577851
$ws.Range("G2").Value = 45.48490528368133
$ws.Range("H2").Value = 17.93919805581172
$ws.Range("I2").Value = 26.52254716885239
$ws.Range("J2").Value = 10.37556985123617
$ws.Range("L2").Value = 12.50625070224778
$ws.Range("M2").Value = 17.54293314813208
$ws.Range("N2").Value = 19.36827236398478

# Row 3
$ws.Range("B3").Value = 16.53086131401737
$ws.Range("D3").Value = 8.258831802346183
$ws.Range("E3").Value = 13.62527057450888
$ws.Range("F3").Value = 37.59297653878034
$ws.Range("G3").Value = 45.06117045509102
$ws.Range("H3").Value = 17.93776617674463
$ws.Range("I3").Value = 26.66993717242959
$ws.Range("J3").Value = 10.30693997531794
$ws.Range("L3").Value = 12.37854090831665
$ws.Range("M3").Value = 17.41233244829314
$ws.Range("N3").Value = 19.42071891221026

# Row 4
$ws.Range("B4").Value = 16.37622787236128
$ws.Range("D4").Value = 8.200702558340438
$ws.Range("E4").Value = 13.50884790899317
$ws.Range("F4").Value = 37.57207838993381
$ws.Range("G4").Value = 44.8159784091018
$ws.Range("H4").Value = 17.94140507822029
$ws.Range("I4").Value = 26.76534760970285
$ws.Range("J4").Value = 10.26409087708235
$ws.Range("L4").Value = 12.30168215020788
$ws.Range("M4").Value = 17.33438754171001
$ws.Range("N4").Value = 19.45493262601634

# Row 5
$ws.Range("B5").Value = 16.3134689353842
$ws.Range("D5").Value = 8.176672179884285
$ws.Range("E5").Value = 13.46077359841006
$ws.Range("F5").Value = 37.56636783107985
$ws.Range("G5").Value = 44.71993324305627
$ws.Range("H5").Value = 17.94402264911372
$ws.Range("I5").Value = 26.80546489280177
$ws.Range("J5").Value = 10.24645628822381
$ws.Range("L5").Value = 12.27078293902196
$ws.Range("M5").Value = 17.30321506575067
$ws.Range("N5").Value = 19.46938188608613

# Row 6
$ws.Range("B6").Value = 16.30306546397866
$ws.Range("D6").Value = 8.172661538532511
$ws.Range("E6").Value = 13.45275353552049
$ws.Range("F6").Value = 37.56558902342875
$ws.Range("G6").Value = 44.70422179662227
$ws.Range("H6").Value = 17.94452577192962
$ws.Range("I6").Value = 26.81220107918368
$ws.Range("J6").Value = 10.2435177706049
$ws.Range("L6").Value = 12.26567846400776
$ws.Range("M6").Value = 17.29807532706604
$ws.Range("N6").Value = 19.47181182341295

# Row 7
$ws.Range("B7").Value = 16.37538035087907
$ws.Range("D7").Value = 8.200379850412343
$ws.Range("E7").Value = 13.50820208319081
$ws.Range("F7").Value = 37.57199001645373
$ws.Range("G7").Value = 44.81466729877246
$ws.Range("H7").Value = 17.94143578786979
$ws.Range("I7").Value = 26.76588363586892
$ws.Range("J7").Value = 10.26385374521963
$ws.Range("L7").Value = 12.30126368732581
$ws.Range("M7").Value = 17.33396471305281
$ws.Range("N7").Value = 19.45512543990625

# Row 8
$ws.Range("B8").Value = 16.69653341465763
$ws.Range("D8").Value = 8.319666930278077
$ws.Range("E8").Value = 13.74727495297786
$ws.Range("F8").Value = 37.62486513789883
$ws.Range("G8").Value = 45.33575763652028
$ws.Range("H8").Value = 17.93776596676181
$ws.Range("I8").Value = 26.57234857335912
$ws.Range("J8").Value = 10.35205683254005
$ws.Range("L8").Value = 12.46191185203251
$ws.Range("M8").Value = 17.49745226514803
$ws.Range("N8").Value = 19.3859391608195

# Row 9
$ws.Range("B9").Value = 17.32878993875579
$ws.Range("D9").Value = 8.542150169640204
$ws.Range("E9").Value = 14.19429245788914
$ws.Range("F9").Value = 37.81688829896886
$ws.Range("G9").Value = 46.47144718025129
$ws.Range("H9").Value = 17.96645303432089
$ws.Range("I9").Value = 26.23173567459775
$ws.Range("J9").Value = 10.51915841022514
$ws.Range("L9").Value = 12.78792610209617
$ws.Range("M9").Value = 17.83466228359238
$ws.Range("N9").Value = 19.2661735627803

# Row 10
$ws.Range("B10").Value = 17.7905522609755
$ws.Range("D10").Value = 8.69816277654245
$ws.Range("E10").Value = 14.50803382731119
$ws.Range("F10").Value = 38.0114259838528
$ws.Range("G10").Value = 47.36793444938738
$ws.Range("H10").Value = 18.009410246683
$ws.Range("I10").Value = 26.00511067449754
$ws.Range("J10").Value = 10.63804997106031
$ws.Range("L10").Value = 13.0322923825122
$ws.Range("M10").Value = 18.09088410939408
$ws.Range("N10").Value = 19.1878110572166

# Row 11
$ws.Range("B11").Value = 17.99914928116917
$ws.Range("D11").Value = 8.767424409012472
$ws.Range("E11").Value = 14.64729686135653
$ws.Range("F11").Value = 38.11139189473069
$ws.Range("G11").Value = 47.78757567774488
$ws.Range("H11").Value = 18.03368949171342
$ws.Range("I11").Value = 25.90712346229374
$ws.Range("J11").Value = 10.69123100113982
$ws.Range("L11").Value = 13.1441145578654
$ws.Range("M11").Value = 18.2089226721334
$ws.Range("N11").Value = 19.15423856230283

# Row 12
$ws.Range("B12").Value = 18.0778596284761
$ws.Range("D12").Value = 8.793398339795486
$ws.Range("E12").Value = 14.6995124217843
$ws.Range("F12").Value = 38.15087914466042
$ws.Range("G12").Value = 47.94803671281062
$ws.Range("H12").Value = 18.04356213403139
$ws.Range("I12").Value = 25.87075128018335
$ws.Range("J12").Value = 10.7112339262313
$ws.Range("L12").Value = 13.18652011841671
$ws.Range("M12").Value = 18.25380280041183
$ws.Range("N12").Value = 19.14182286222885

# Row 13
$ws.Range("B13").Value = 18.0609215334428
$ws.Range("D13").Value = 8.787815810315998
$ws.Range("E13").Value = 14.68829036179025
$ws.Range("F13").Value = 38.14230256029521
$ws.Range("G13").Value = 47.91341191947339
$ws.Range("H13").Value = 18.04140575285368
$ws.Range("I13").Value = 25.87855207077853
$ws.Range("J13").Value = 10.7069320580386
$ws.Range("L13").Value = 13.17738514671868
$ws.Range("M13").Value = 18.24412949492239
$ws.Range("N13").Value = 19.14448358928464

# Row 14
$ws.Range("B14").Value = 18.00563085226976
$ws.Range("D14").Value = 8.76956643511676
$ws.Range("E14").Value = 14.65160322250445
$ws.Range("F14").Value = 38.11460792734118
$ws.Range("G14").Value = 47.80074665752141
$ws.Range("H14").Value = 18.03448813746111
$ws.Range("I14").Value = 25.90411640546961
$ws.Range("J14").Value = 10.69287940330659
$ws.Range("L14").Value = 13.14760226008148
$ws.Range("M14").Value = 18.21261153162333
$ws.Range("N14").Value = 19.15321115845358

# Row 15
$ws.Range("B15").Value = 17.97172508996147
$ws.Range("D15").Value = 8.758354857256306
$ws.Range("E15").Value = 14.62906286672711
$ws.Range("F15").Value = 38.09785619111543
$ws.Range("G15").Value = 47.73193360333558
$ws.Range("H15").Value = 18.03033918201132
$ws.Range("I15").Value = 25.91987079500162
$ws.Range("J15").Value = 10.68425392286585
$ws.Range("L15").Value = 13.12936631965753
$ws.Range("M15").Value = 18.19332856313653
$ws.Range("N15").Value = 19.15859576009452

# Row 16
$ws.Range("B16").Value = 17.77688444931362
$ws.Range("D16").Value = 8.693601267665208
$ws.Range("E16").Value = 14.49886107043505
$ws.Range("F16").Value = 38.00512238587882
$ws.Range("G16").Value = 47.34073522875143
$ws.Range("H16").Value = 18.00791866699815
$ws.Range("I16").Value = 26.01161704272715
$ws.Range("J16").Value = 10.63455580160519
$ws.Range("L16").Value = 13.02499489966264
$ws.Range("M16").Value = 18.08319713349371
$ws.Range("N16").Value = 19.19004677491019

# Row 17
$ws.Range("B17").Value = 17.65693067363656
$ws.Range("D17").Value = 8.653433362821266
$ws.Range("E17").Value = 14.41808385461622
$ws.Range("F17").Value = 37.95115964941117
$ws.Range("G17").Value = 47.10367160470723
$ws.Range("H17").Value = 17.99537621760916
$ws.Range("I17").Value = 26.0692075174494
$ws.Range("J17").Value = 10.60383219088784
$ws.Range("L17").Value = 12.96111099249963
$ws.Range("M17").Value = 18.01599219572167
$ws.Range("N17").Value = 19.20987175793717

# Row 18
$ws.Range("B18").Value = 17.5878029919137
$ws.Range("D18").Value = 8.630169431347175
$ws.Range("E18").Value = 14.37129861915059
$ws.Range("F18").Value = 37.92120249201211
$ws.Range("G18").Value = 46.96844096387509
$ws.Range("H18").Value = 17.98860828555306
$ws.Range("I18").Value = 26.10281264047707
$ws.Range("J18").Value = 10.58607609277945
$ws.Range("L18").Value = 12.92443066966655
$ws.Range("M18").Value = 17.97747909007292
$ws.Range("N18").Value = 19.22146992558961

# Row 19
$ws.Range("B19").Value = 17.56437682637731
$ws.Range("D19").Value = 8.622265365314931
$ws.Range("E19").Value = 14.35540295301718
$ws.Range("F19").Value = 37.91124561152969
$ws.Range("G19").Value = 46.92285138756924
$ws.Range("H19").Value = 17.98639347578924
$ws.Range("I19").Value = 26.11427332597083
$ws.Range("J19").Value = 10.58004981139627
$ws.Range("L19").Value = 12.91202337118631
$ws.Range("M19").Value = 17.96446445994302
$ws.Range("N19").Value = 19.22543044747353

# Row 20
$ws.Range("B20").Value = 17.66971428960499
$ws.Range("D20").Value = 8.65772596703294
$ws.Range("E20").Value = 14.42671646131316
$ws.Range("F20").Value = 37.95679233288013
$ws.Range("G20").Value = 47.12879236737304
$ws.Range("H20").Value = 17.99666522465241
$ws.Range("I20").Value = 26.06302717808281
$ws.Range("J20").Value = 10.6071115780237
$ws.Range("L20").Value = 12.96790516447349
$ws.Range("M20").Value = 18.02313187065858
$ws.Range("N20").Value = 19.20774114216728

# Row 21
$ws.Range("B21").Value = 18.02187923707956
$ws.Range("D21").Value = 8.774933674636417
$ws.Range("E21").Value = 14.66239342345456
$ws.Range("F21").Value = 38.12269836045456
$ws.Range("G21").Value = 47.83379824930921
$ws.Range("H21").Value = 18.03650161638575
$ws.Range("I21").Value = 25.89658764050199
$ws.Range("J21").Value = 10.69701073312978
$ws.Range("L21").Value = 13.15634882983817
$ws.Range("M21").Value = 18.2218644486038
$ws.Range("N21").Value = 19.15063959355022

# Row 22
$ws.Range("B22").Value = 18.25037173573821
$ws.Range("D22").Value = 8.850051907871261
$ws.Range("E22").Value = 14.81337917134987
$ws.Range("F22").Value = 38.24063221061692
$ws.Range("G22").Value = 48.30353637585419
$ws.Range("H22").Value = 18.06649072713049
$ws.Range("I22").Value = 25.79208460528974
$ws.Range("J22").Value = 10.75497219094687
$ws.Range("L22").Value = 13.27984740571888
$ws.Range("M22").Value = 18.35279078819172
$ws.Range("N22").Value = 19.11505393951012

# Row 23
$ws.Range("B23").Value = 18.12859623139618
$ws.Range("D23").Value = 8.810098281228329
$ws.Range("E23").Value = 14.73308094497749
$ws.Range("F23").Value = 38.17682533489498
$ws.Range("G23").Value = 48.05205657292347
$ws.Range("H23").Value = 18.050124290664
$ws.Range("I23").Value = 25.84746891417258
$ws.Range("J23").Value = 10.72411147535934
$ws.Range("L23").Value = 13.21391387233589
$ws.Range("M23").Value = 18.28282796511532
$ws.Range("N23").Value = 19.13388835010226

# Row 24
$ws.Range("B24").Value = 17.66393532426023
$ws.Range("D24").Value = 8.655785812311839
$ws.Range("E24").Value = 14.42281473380205
$ws.Range("F24").Value = 37.95424247261219
$ws.Range("G24").Value = 47.11743195630623
$ws.Range("H24").Value = 17.99608108478332
$ws.Range("I24").Value = 26.06581976571147
$ws.Range("J24").Value = 10.6056292553896
$ws.Range("L24").Value = 12.96483336912689
$ws.Range("M24").Value = 18.01990363523135
$ws.Range("N24").Value = 19.20870376887816

# Row 25
$ws.Range("B25").Value = 17.15792886925303
$ws.Range("D25").Value = 8.483236712919377
$ws.Range("E25").Value = 14.07584327693268
$ws.Range("F25").Value = 37.7555134434437
$ws.Range("G25").Value = 46.15275432677034
$ws.Range("H25").Value = 17.95484731071731
$ws.Range("I25").Value = 26.31972532628962
$ws.Range("J25").Value = 10.47461215058003
$ws.Range("L25").Value = 12.69874479545073
$ws.Range("M25").Value = 17.74182783894807
$ws.Range("N25").Value = 19.29687743837015

Write-Output "Updated loading_percent values for rows 2-25 (380 kV case)."
